$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.41"
$ws.Range("E2").Value = "'-2.90%"
$ws.Range("D3").Value = "'35.30"
$ws.Range("E3").Value = "'-0.48%"
$ws.Range("D4").Value = "'5.077"
$ws.Range("E4").Value = "'-0.13%"
$ws.Range("D5").Value = "'0.07915"
$ws.Range("E5").Value = "'-3.00%"
$ws.Range("D6").Value = "'1.893"
$ws.Range("E6").Value = "'-8.52%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.051"
$ws.Range("E7").Value = "'-1.86%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'7.783"
$ws.Range("E8").Value = "'-2.03%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9302"
$ws.Range("E9").Value = "'0.41%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1375"
$ws.Range("E10").Value = "'30.70%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1894"
$ws.Range("E11").Value = "'-1.55%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09167"
$ws.Range("E12").Value = "'-0.62%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03439"
$ws.Range("E13").Value = "'-5.81%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09843"
$ws.Range("E14").Value = "'-0.47%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001404"
$ws.Range("E15").Value = "'-1.61%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005820"
$ws.Range("E16").Value = "'1.06%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.531"
$ws.Range("E17").Value = "'1.53%"
$ws.Range("D18").Value = "'2.986"
$ws.Range("E18").Value = "'4.47%"
$ws.Range("D19").Value = "'0.3411"
$ws.Range("E19").Value = "'0.55%"
$ws.Range("D20").Value = "'0.1305"
$ws.Range("E20").Value = "'0.34%"
$ws.Range("D21").Value = "'5.046"
$ws.Range("E21").Value = "'-1.04%"
$ws.Range("D22").Value = "'0.2403"
$ws.Range("E22").Value = "'8.61%"
$ws.Range("D23").Value = "'0.04508"
$ws.Range("E23").Value = "'-0.80%"
$ws.Range("D24").Value = "'0.001215"
$ws.Range("E24").Value = "'-0.93%"
$ws.Range("D25").Value = "'0.004762"
$ws.Range("E25").Value = "'-0.37%"
$ws.Range("D26").Value = "'0.0001232"
$ws.Range("E26").Value = "'-1.57%"
$ws.Range("D27").Value = "'0.0003006"
$ws.Range("E27").Value = "'-32.49%"
$ws.Range("D39").Value = "'0.01852"
$ws.Range("E39").Value = "'-5.98%"
$ws.Range("D40").Value = "'0.04753"
$ws.Range("E40").Value = "'-2.81%"
$ws.Range("D41").Value = "'0.007306"
$ws.Range("E41").Value = "'-3.48%"
$ws.Range("D42").Value = "'0.009625"
$ws.Range("E42").Value = "'7.30%"
$ws.Range("E43").Value = "'-4.21%"
$ws.Range("E44").Value = "'-1.71%"
$ws.Range("D45").Value = "'0.01095"
$ws.Range("E45").Value = "'-5.83%"
$ws.Range("D46").Value = "'0.00006261"
$ws.Range("E46").Value = "'-5.43%"
$ws.Range("E47").Value = "'0.08%"
$ws.Range("D48").Value = "'64.66"
$ws.Range("E48").Value = "'-64.78%"
$ws.Range("E49").Value = "'10.52%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.08%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.08%"
